$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Expand the table from 1 column (A1:A15) to 3 columns (A1:C15) ---------
$lo.Resize($ws.Range("A1:C15")) | Out-Null

# Stamp the same cell style that column A already carries onto the newly
# created B/C columns, row by row (header row keeps the header style s="2",
# every data row keeps the data style s="1"), so every row ends up with a
# real (if empty) cell per column, matching the original table's per-row
# shape/formatting.
for ($r = 1; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Copy($ws.Cells.Item($r, 2)) | Out-Null
    $ws.Cells.Item($r, 1).Copy($ws.Cells.Item($r, 3)) | Out-Null
}

# --- Rename columns via the header cells: tipo | prefijo | codigo ---------
$ws.Range("A1").Value = "tipo"
$ws.Range("B1").Value = "prefijo"
$ws.Range("C1").Value = "codigo"

# --- Replace the old "comprobantes a desestimar" numeric listing with the
#     new single data row, clearing every other row's contents (keeping the
#     now-stamped styles intact). ------------------------------------------
$ws.Range("A2:C15").ClearContents() | Out-Null

$ws.Range("A2").Value = "FACB2"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 30970

$ws.Range("C3").Select() | Out-Null

# The original column (now "codigo", column C) keeps its custom width;
# column A no longer needs the custom width it used to carry.
$ws.Columns.Item(3).ColumnWidth = 21.17
$ws.Columns.Item(1).UseStandardWidth = $true

$ws.PageSetup.Orientation = 1

Write-Output "done"
